$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per column letter (applied identically to rows 2 and 3)
$values = @{
    "D"  = 0.0155
    "E"  = 0.005180000000000001
    "G"  = 1.263157894736842
    "H"  = 1.263157894736842
    "I"  = 0.9672830725462304
    "J"  = 0.9035726009410219
    "K"  = 5.83
    "L"  = 0.829302987197724
    "M"  = 5.87
    "N"  = 0.07402269861286255
    "O"  = 1.006861063464837
    "P"  = 5.87
    "Q"  = 0.07402269861286255
    "R"  = 1.006861063464837
    "U"  = 0.35
    "V"  = 0.004413619167717528
    "W"  = 0.239917695473251
    "X"  = 0.06454268619014099
    "Y"  = 0.1753750092831101
    "Z"  = 0.2914956254923913
    "AA" = 0.26338746048909
    "AB" = 0.06454268619014099
    "AC" = 0.198844774298949
    "AG" = -0.35
    "AJ" = -0.004433185560481317
    "AK" = -0.01455301455301455
    "AM" = -0.005
    "AP" = -0.05131964809384164
    "AQ" = -1360
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
